$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.014.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.98%  '
$ws.Range("D3").Value = '''1.895.21'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.42%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '''249.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = '''0.4980'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("D8").Value = '''44.85'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.64%  '
$ws.Range("D9").Value = '''0.2960'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.95%  '
$ws.Range("D10").Value = '''0.06656'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("D11").Value = '''1.898.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.02%  '
$ws.Range("D12").Value = '''16.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("D13").Value = '''0.07266'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("D14").Value = '''0.6772'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.65%  '
$ws.Range("D15").Value = '''85.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").Value = '''4.853'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("D17").Value = '''30.027.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.97%  '
$ws.Range("D18").Value = '''0.000008053'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +10.31%  '
$ws.Range("D19").Value = '''0.9997'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '''12.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.47%  '
$ws.Range("D21").Value = '''2.149.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = '''1.008'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '''4.760'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.72%  '
$ws.Range("D24").Value = '''9.238'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.41%  '
$ws.Range("D25").Value = '''5.646'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.67%  '
$ws.Range("D26").Value = '''147.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").Value = '''131.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").Value = '''16.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").Value = '''1.958'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.20%  '
$ws.Range("D30").Value = '''1.392'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").Value = '''4.226'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.62%  '
$ws.Range("D32").Value = '''0.08781'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.37%  '
$ws.Range("D33").Value = '''3.934'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.32%  '
$ws.Range("D34").Value = '''0.05114'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("D35").Value = '''1.124'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.95%  '
$ws.Range("D36").Value = '''0.7048'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.43%  '
$ws.Range("D37").Value = '''2.698'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '''2.792'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("D39").Value = '''2.238'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("D40").Value = '''0.9577'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = '''0.01660'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.03%  '
$ws.Range("D42").Value = '''6.025'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.4228'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.03%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''103.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("D46").Value = '''7.469'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").Value = '''0.1261'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.56%  '
$ws.Range("D48").Value = '''0.05767'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.56%  '
$ws.Range("D49").Value = '''32.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.13%  '
$ws.Range("D50").Value = '''8.316'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("D51").Value = '''0.3737'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.68%  '
